$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of column J (|S*|/n average)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Rows 14-17: summary labels and formulas
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Build the bold/size-12/vertical-centered format once on a helper cell,
# then copy its format onto B14:B17 so only one new style is created.
$helper = $ws.Range("Z100")
$helper.Font.Bold = $true
$helper.Font.Size = 12
$helper.VerticalAlignment = -4108

$target = $ws.Range("B14:B17")
$helper.Copy()
$target.PasteSpecial(-4122)  # xlPasteFormats
$helper.Clear()
$excel.CutCopyMode = 0

# Selection as in the saved file
$ws.Range("A14:B17").Select() | Out-Null

# Print setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
